$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "70.122.73"
$ws.Range("E2").Value = "  -2.97%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.519.57"
$ws.Range("E3").Value = "  -4.47%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.26"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "167.92"
$ws.Range("E6").Value = "  -4.07%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.517.75"
$ws.Range("E9").Value = "  -4.46%  "
$ws.Range("E10").Value = "  -5.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.348"
$ws.Range("E12").Value = "  -3.05%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.90"
$ws.Range("E13").Value = "  -0.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.980.04"
$ws.Range("E14").Value = "  -4.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "70.048.53"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("E16").Value = "  -5.39%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.16"
$ws.Range("E17").Value = "  -2.28%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.525.79"
$ws.Range("E18").Value = "  -3.72%  "
$ws.Range("E19").Value = "  -2.10%  "
$ws.Range("E20").Value = "  -6.03%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "351.70"
$ws.Range("E21").Value = "  -6.25%  "
$ws.Range("E22").Value = "  -3.87%  "
$ws.Range("E23").Value = "  -3.46%  "
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "69.47"
$ws.Range("E25").Value = "  -2.79%  "
$ws.Range("E26").Value = "  -5.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.01"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("E28").Value = "  -4.46%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -3.99%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.91"
$ws.Range("E31").Value = "  -0.29%  "
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "465.87"
$ws.Range("E33").Value = "  -5.08%  "
$ws.Range("E34").Value = "  -2.12%  "
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("E36").Value = "  +3.20%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "153.34"
$ws.Range("E37").Value = "  -4.22%  "
$ws.Range("E38").Value = "  +0.70%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.48"
$ws.Range("E39").Value = "  -3.55%  "
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  -1.84%  "
$ws.Range("E42").Value = "  -1.19%  "
$ws.Range("E43").Value = "  -6.78%  "
$ws.Range("E44").Value = "  -13.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.32"
$ws.Range("E45").Value = "  -9.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "38.20"
$ws.Range("E46").Value = "  -2.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "143.13"
$ws.Range("E47").Value = "  -4.70%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.531"
$ws.Range("E48").Value = "  -2.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.50"
$ws.Range("E49").Value = "  -3.25%  "
$ws.Range("E50").Value = "  -4.41%  "
$ws.Range("E51").Value = "  -1.18%  "
